$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ireland-manual")

# New daily COVID rows appended to the Ireland-manual sheet (rows 194-197):
# same Lat/Long as the rest of the sheet, two new dates (2020-04-27,
# 2020-04-28), each with a "confirmed" and a "death" entry.
$newRows = @(
    @{ Row = 194; Date = "2020-04-27"; Cases = 386; Type = "confirmed" },
    @{ Row = 195; Date = "2020-04-27"; Cases = 18;  Type = "death" },
    @{ Row = 196; Date = "2020-04-28"; Cases = 229; Type = "confirmed" },
    @{ Row = 197; Date = "2020-04-28"; Cases = 59;  Type = "death" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = "Ireland"
    $ws.Cells.Item($row, 3).Value = 53.1424
    $ws.Cells.Item($row, 4).Value = -7.6921

    # Force the date column to be stored as text (matches the rest of the
    # sheet, where dates are plain text strings, not date serials).
    $dateCell = $ws.Cells.Item($row, 5)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date

    $ws.Cells.Item($row, 6).Value = $r.Cases
    $ws.Cells.Item($row, 7).Value = $r.Type
}

# Update the active selection to match the end of the newly appended data.
$ws.Range("F198").Select()
